# SCD0012-003 - Admin SLN melakukan meng-aktifkan atau me-nonaktifkan parameter KPI
# Commit: "Update Excel SCD0011 until SCD0016"
#
# Changes applied:
#  1. Rename the worksheet from "SCD0213" to "SCD0012".
#  2. Update the TC_ID cell (B2) from "DGS-228" to "SCD0012-003".
#  3. Widen column B to fit the new, longer TC_ID text.
#  4. Reset the view: scroll back to the top-left (A1) and select B2
#     instead of the previous T2 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "SCD0012"

# 2. Update TC_ID value in B2 (shared-string table is re-packed automatically)
$ws.Range("B2").Value = "SCD0012-003"

# 3. Re-fit column B so the new (longer) TC_ID value is fully visible
$ws.Columns.Item(2).ColumnWidth = 11.6

# 4. Reset scroll position to A1 and move the selection to B2
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
[void]$ws.Range("B2").Select()
